# Apply the "BBox transforms in detection file" edit to the workbook.
#
# Summary of the change (from the OOXML diff):
#  - On the `dataset_params` sheet, two new configuration rows are inserted
#    right before the existing "plan" row (old row 13):
#       row 13: var_name = "cache_rate", manual_value (col B) = 0.3
#       row 14: var_name = "ds_type"
#    The old row 13 ("plan" / tune_p / manual_p / E=3 formatting row) shifts
#    down to become row 15.
#  - The `dataset_params` sheet becomes the active sheet / active tab of the
#    workbook, with the active selection on cell E14 (next to the new
#    "ds_type" row).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("dataset_params")

# Insert two fresh rows above the old row 13 ("plan" row), pushing it (and
# everything below it) down to rows 15+.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(14).Insert()

# New row 13: cache_rate
$ws.Range("A13").Value = "cache_rate"
$ws.Range("B13").Value = 0.3

# New row 14: ds_type
$ws.Range("A14").Value = "ds_type"

# Make dataset_params the active sheet/tab, with E14 selected (this also
# normalizes the other sheets' lingering "B:B"-style multi-area selections
# left over from the previous editor, and clears tabSelected on whichever
# sheet used to be active).
$ws.Activate()
$ws.Range("E14").Select()
